$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new quarter-end dates in R1:T1, matching the style of the existing header cells ---
$ws.Range("Q1").Copy()
$ws.Range("R1:T1").PasteSpecial(-4122)
$ws.Range("R1").Value = "31/12/2023"
$ws.Range("S1").Value = "31/03/2024"
$ws.Range("T1").Value = "30/06/2024"

# --- Data rows: balance-sheet / income-statement figures for the three new quarters ---
$rsData = @{
    2 = @(682102.0159999999, 696284.992, 706393.9840000001)
    3 = @(502704.992, 507967.008, 500148.992)
    4 = @(54477, 82179, 66988)
    5 = @(0, 0, 0)
    6 = @(196536, 148743.008, 152948.992)
    7 = @(232856, 253506, 260536.992)
    8 = @(0, 0, 0)
    9 = @(10336, 10915, 4359)
    10 = @(8500, 12624, 15316)
    11 = @(0, 0, 0)
    12 = @(12101, 11998, 12470)
    13 = @(0, 0, 0)
    14 = @(0, 0, 0)
    15 = @(0, 0, 0)
    16 = @(0, 0, 0)
    17 = @(0, 0, 0)
    18 = @(0, 0, 0)
    19 = @(6400, 6178, 5911)
    20 = @(0, 0, 0)
    21 = @(0, 0, 0)
    22 = @(0, 0, 0)
    23 = @(146610, 151924.992, 169272.992)
    24 = @(20686, 24395, 24502)
    25 = @(0, 0, 0)
    26 = @(682102.0159999999, 696284.992, 706393.9840000001)
    27 = @(177724, 169148.992, 154610)
    28 = @(34652, 32214, 29140)
    29 = @(59198, 64033, 69998)
    30 = @(0, 0, 0)
    31 = @(17682, 16465, 16801)
    32 = @(0, 0, 0)
    33 = @(27473, 33113, 12327)
    34 = @(38719, 23324, 26344)
    35 = @(0, 0, 0)
    36 = @(0, 0, 0)
    37 = @(90745, 92342, 103738)
    38 = @(84245, 86713, 98775)
    39 = @(0, 0, 0)
    40 = @(294, 0, 0)
    41 = @(0, 0, 0)
    42 = @(0, 0, 0)
    43 = @(6206, 5629, 4963)
    44 = @(0, 0, 0)
    45 = @(0, 0, 0)
    46 = @(0, 0, 0)
    47 = @(413632.992, 434793.984, 448046.016)
    48 = @(192392, 192392, 336148)
    49 = @(-16141, -15123, -21650)
    50 = @(0, 0, 0)
    51 = @(235516, 235516, 92127)
    52 = @(0, 20142, 39552)
    53 = @(0, 0, 0)
    54 = @(0, 0, 0)
    55 = @(1866, 1867, 1869)
    56 = @(0, 0, 0)
    59 = @(221946.976, 162131.008, 192002)
    60 = @(-94101, -66695, -84783)
    61 = @(127846.016, 95436, 107219)
    62 = @(-45842, -31567, -35895)
    63 = @(-35899, -29861, -35080)
    64 = @(0, 0, 0)
    65 = @(1013, 0, -339)
    66 = @(50, 307, 0)
    67 = @(0, 0, 0)
    68 = @(-3192, -3534, -3937)
    69 = @(1983, 2547, 2313)
    70 = @(-5175, -6081, -6250)
    74 = @(43976, 30781, 31968)
    75 = @(-7165, -4174, -5456)
    76 = @(-1483, -222, -267)
    79 = @(0, 0, 0)
    80 = @(35328, 26385, 26245)
}

foreach ($row in $rsData.Keys) {
    $vals = $rsData[$row]
    $ws.Cells.Item($row, 18).Value = $vals[0]
    $ws.Cells.Item($row, 19).Value = $vals[1]
    $ws.Cells.Item($row, 20).Value = $vals[2]
}
